$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MaxCapacityGroup")

# Update the note text used by the PumpedHydro rows (shared string index 21,
# currently "GW") in place so every cell referencing it picks up the new text.
$ws.Cells.Replace("GW", "GW- Assume no new pumped hydro is added above the already installed capacity")

# Update the max capacity values (GW) for PumpedHydro rows 2:8
$ws.Range("C2:C8").Value = 0.095

# Update the max capacity values (GW) and notes for HydroGeneration rows 9:15
$ws.Range("C9:C15").Value = 2.058
$ws.Range("D9:D15").Value = "GW- Assume no new hydro is added above the already installed capacity"

# Update the active selection to match the saved view state
$ws.Range("F17").Select()

$wb.Save()
